$wb = $excel.ActiveWorkbook

# Update "想去人数" (interested-count) figures on both the "展览" sheet
# and the "全部类型" sheet, which mirror the same rows.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 6168
    $ws.Range("F5").Value = 1002
    $ws.Range("F6").Value = 101
}
